# Update Work Week and Social Spending
# (commit message is generic; the actual edit refreshes the Mali GDP-per-Capita
# series in the "Data" sheet and appends years 2011-2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Updated values for existing rows 2..62 (years 1950..2010), column E ---
$updatedValues = @("728","740","752","765","778","791","803","816","830","845","851","840","827","866","888","881","899","910","923","899","934","945","971","926","883","982","1090","1148","1101","1345","1168","1111","1030","1057","1084","1089","1138","1138","1116","1205","1191","1186.75768210614","1115.40597097193","1115.10017260117","1111.78536659544","1100.2432497145","1130.50966495942","1160.95166520433","1167.43257807046","1192.31681889464","1153.04614535542","1309.25583970687","1321.33882657338","1391.26347128294","1366.98306657794","1414.79273856253","1438.20793270406","1441.20593670952","1460.49598003884","1477.41837438766","1504.99561226943")

$startRow = 2
$endRow = 62

# The Data column stores numeric-looking readings as text (matches the
# source workbook's convention), so force text formatting before writing.
$ws.Range("E$startRow`:E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $updatedValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $updatedValues[$i]
}

# --- New rows 63..68 for years 2011..2016 ---
$newYears = @(2011, 2012, 2013, 2014, 2015, 2016)

# Values per the refreshed series for 2011-2016
$newRowValues = @("1500", "1448", "1443", "1505", "1553", "1594")

$countryCode = 466
$countryName = "Mali"
$indicator = "GDP per Capita"

$newStartRow = 63
$ws.Range("E$newStartRow`:E68").NumberFormat = "@"

for ($i = 0; $i -lt $newYears.Count; $i++) {
    $row = $newStartRow + $i
    $ws.Cells.Item($row, 1).Value = $countryCode
    $ws.Cells.Item($row, 2).Value = $countryName
    $ws.Cells.Item($row, 3).Value = $indicator
    $ws.Cells.Item($row, 4).Value = $newYears[$i]
    $ws.Cells.Item($row, 5).Value = $newRowValues[$i]
}
